$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.045.82"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "3.518.75"
$ws.Range("E3").Value = "  -3.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "3.512.72"
$ws.Range("E8").Value = "  -3.27%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.578"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").Value = "4.085.43"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "614.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.12%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.532.44"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "69.100.16"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.885"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -6.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.44%  "
$ws.Range("E30").Value = "  -6.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.54%  "
$ws.Range("E32").Value = "  -5.85%  "
$ws.Range("E33").Value = "  -5.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "610.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -12.91%  "
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("D43").Value = "3.396.85"
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("D46").Value = "0.0₃0697"
$ws.Range("E46").Value = "  -4.97%  "
$ws.Range("E47").Value = "  -6.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.80%  "
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.55%  "
